$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column before the existing "Actioned Yes/No" column (E),
#    shifting Actioned Yes/No -> F, Comments -> G, Notes -> H.
$ws.Columns("E:E").Insert()

# 2. New column header + priority values
$ws.Range("H11").Value = "Just need to push this update across the other tabs, when time permits"
$ws.Range("E2").Value = "Priority"

$ws.Range("E3").Value = 1
$ws.Range("E4").Value = 2
$ws.Range("E5").Value = 1
$ws.Range("E6").Value = 3
$ws.Range("E7").Value = 2
$ws.Range("E8").Value = 2
$ws.Range("E9").Value = 1
$ws.Range("E10").Value = 2

# Row 11 also picks up a "Yes" in the Follow Up Questions actioned column
# and a priority value (the note text was already set above)
$ws.Range("D11").Value = "Yes"
$ws.Range("E11").Value = 1

$ws.Range("E14").Value = 1
$ws.Range("E15").Value = 1
$ws.Range("E16").Value = 1
$ws.Range("E17").Value = 1
$ws.Range("E18").Value = 1

$ws.Range("E20").Value = 3
$ws.Range("E21").Value = 3
$ws.Range("E22").Value = 3

$ws.Range("E25").Value = 2
$ws.Range("E26").Value = 2
$ws.Range("E27").Value = 2
$ws.Range("E28").Value = 2
$ws.Range("E29").Value = 2
